$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '''36.946.23'
$ws.Cells.Item(2, 5).Value = '  -0.50%  '
$ws.Cells.Item(3, 4).Value = '''2.044.53'
$ws.Cells.Item(3, 5).Value = '  -0.18%  '
$ws.Cells.Item(4, 5).Value = '  -0.07%  '
$ws.Cells.Item(5, 4).Value = '''247.98'
$ws.Cells.Item(5, 5).Value = '  -0.67%  '
$ws.Cells.Item(6, 4).Value = '''0.661'
$ws.Cells.Item(6, 5).Value = '  +0.73%  '
$ws.Cells.Item(7, 5).Value = '  +0.01%  '
$ws.Cells.Item(8, 4).Value = '''56.03'
$ws.Cells.Item(8, 5).Value = '  +1.61%  '
$ws.Cells.Item(9, 4).Value = '''0.379'
$ws.Cells.Item(9, 5).Value = '  -0.14%  '
$ws.Cells.Item(10, 4).Value = '''0.0779'
$ws.Cells.Item(10, 5).Value = '  +2.97%  '
$ws.Cells.Item(11, 5).Value = '  +1.24%  '
$ws.Cells.Item(12, 4).Value = '''15.71'
$ws.Cells.Item(12, 5).Value = '  +3.66%  '
$ws.Cells.Item(13, 4).Value = '''2.338.49'
$ws.Cells.Item(13, 5).Value = '  -0.25%  '
$ws.Cells.Item(14, 4).Value = '''5.59'
$ws.Cells.Item(14, 5).Value = '  +6.06%  '
$ws.Cells.Item(15, 4).Value = '''0.790'
$ws.Cells.Item(15, 5).Value = '  -4.44%  '
$ws.Cells.Item(16, 4).Value = '''2.044.70'
$ws.Cells.Item(16, 5).Value = '  -0.14%  '
$ws.Cells.Item(17, 4).Value = '''36.916.33'
$ws.Cells.Item(17, 5).Value = '  -0.34%  '
$ws.Cells.Item(18, 4).Value = '''16.38'
$ws.Cells.Item(18, 5).Value = '  +13.60%  '
$ws.Cells.Item(19, 4).Value = '''73.79'
$ws.Cells.Item(19, 5).Value = '  +1.82%  '
$ws.Cells.Item(20, 4).Value = '''0.0₃0893'
$ws.Cells.Item(20, 5).Value = '  +1.59%  '
$ws.Cells.Item(21, 4).Value = '''5.30'
$ws.Cells.Item(21, 5).Value = '  +0.25%  '
$ws.Cells.Item(22, 4).Value = '''235.61'
$ws.Cells.Item(22, 5).Value = '  -1.15%  '
$ws.Cells.Item(23, 5).Value = '  +0.10%  '
$ws.Cells.Item(24, 4).Value = '''2.36'
$ws.Cells.Item(24, 5).Value = '  -2.94%  '
$ws.Cells.Item(25, 5).Value = '  +8.10%  '
$ws.Cells.Item(26, 4).Value = '''167.49'
$ws.Cells.Item(26, 5).Value = '  -1.68%  '
$ws.Cells.Item(27, 4).Value = '''9.05'
$ws.Cells.Item(27, 5).Value = '  -1.34%  '
$ws.Cells.Item(28, 4).Value = '''19.73'
$ws.Cells.Item(28, 5).Value = '  -2.96%  '
$ws.Cells.Item(29, 5).Value = '  +0.73%  '
$ws.Cells.Item(30, 5).Value = '  +3.44%  '
$ws.Cells.Item(31, 4).Value = '''4.64'
$ws.Cells.Item(31, 5).Value = '  +0.75%  '
$ws.Cells.Item(32, 4).Value = '''0.0611'
$ws.Cells.Item(32, 5).Value = '  -2.74%  '
$ws.Cells.Item(33, 4).Value = '''4.40'
$ws.Cells.Item(33, 5).Value = '  +0.13%  '
$ws.Cells.Item(34, 5).Value = '  -0.19%  '
$ws.Cells.Item(35, 4).Value = '''0.0870'
$ws.Cells.Item(35, 5).Value = '  +3.64%  '
$ws.Cells.Item(36, 4).Value = '''2.21'
$ws.Cells.Item(36, 5).Value = '  -3.41%  '
$ws.Cells.Item(37, 5).Value = '  -1.34%  '
$ws.Cells.Item(38, 2).Value = 'Cronos'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Cells.Item(38, 4).Value = '''0.106'
$ws.Cells.Item(38, 5).Value = '  -4.80%  '
$ws.Cells.Item(39, 2).Value = 'TrustWalletToken'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Cells.Item(39, 4).Value = '''1.33'
$ws.Cells.Item(39, 5).Value = '  -1.13%  '
$ws.Cells.Item(40, 2).Value = 'HuobiToken'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Cells.Item(40, 4).Value = '''3.18'
$ws.Cells.Item(40, 5).Value = '  +13.87%  '
$ws.Cells.Item(41, 2).Value = 'VeChain'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(41, 4).Value = '''0.0220'
$ws.Cells.Item(41, 5).Value = '  -2.83%  '
$ws.Cells.Item(42, 2).Value = 'THORChain'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Cells.Item(42, 4).Value = '''4.87'
$ws.Cells.Item(42, 5).Value = '  +23.97%  '
$ws.Cells.Item(43, 4).Value = '''17.21'
$ws.Cells.Item(43, 5).Value = '  -5.04%  '
$ws.Cells.Item(44, 5).Value = '  -3.10%  '
$ws.Cells.Item(45, 4).Value = '''95.18'
$ws.Cells.Item(45, 5).Value = '  -2.30%  '
$ws.Cells.Item(46, 4).Value = '''2.42'
$ws.Cells.Item(46, 5).Value = '  +1.86%  '
$ws.Cells.Item(47, 4).Value = '''1.275.87'
$ws.Cells.Item(47, 5).Value = '  -2.21%  '
$ws.Cells.Item(48, 5).Value = '  -2.46%  '
$ws.Cells.Item(49, 4).Value = '''2.229.46'
$ws.Cells.Item(49, 5).Value = '  -0.42%  '
$ws.Cells.Item(50, 4).Value = '''6.66'
$ws.Cells.Item(50, 5).Value = '  -3.33%  '
$ws.Cells.Item(51, 2).Value = 'MultiversX'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Cells.Item(51, 4).Value = '''41.95'
$ws.Cells.Item(51, 5).Value = '  -6.83%  '
